# Automatic data update: append a new day's row (2020-05-23) to the Covid-19
# data table on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data lives inside an Excel Table ("Tabela1"); adding a ListRow grows the
# table's range (and its AutoFilter range) by one row and keeps the table
# structure/formatting machinery consistent.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Copy the formatting from a representative earlier data row (row 64) down
# onto the freshly added row 74 so the new cells pick up the same number
# formats/fonts used elsewhere in the table.
$ws.Range("A64:J64").Copy()
$ws.Range("A74:J74").PasteSpecial(-4122)

# Fill in the new day's values.
$ws.Range("A74").Value = 43974
$ws.Range("B74").Value = 74760
$ws.Range("C74").Value = 341
$ws.Range("D74").Value = 1468
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 18
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 107
$ws.Range("J74").Value = 0

# Match the author's final selection (the newly added row).
$ws.Range("A74:J74").Select() | Out-Null
